$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells keep their original text representation
# (numeric-looking strings and percentages stored as text, not numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "313.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.68%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.136"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.11%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08153"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.40%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.502"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.19%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.964"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.08%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.307"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.07%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9408"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.49%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1967"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.30%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03498"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.53%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09718"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.97%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001409"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.41%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006144"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.35%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.577"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-8.50%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.49%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3465"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.33%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1301"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.88%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.015"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.40%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2494"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.72%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04380"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.09%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001248"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.44%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004730"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "199.50%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.68%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02210"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.99%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05227"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.62%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007565"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.90%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01033"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1394"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.32%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.38%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006827"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "17.31%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
